$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.980461359024048
$ws.Range("B1").Value = 5.39075756072998
$ws.Range("C1").Value = 2.803067684173584
$ws.Range("D1").Value = 2.324413061141968
$ws.Range("E1").Value = 2.107751131057739
